$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "This is git " + "test" (two separate runs, split by the old
#    _GoBack bookmark) become a single run "This is git test".
#    A Find/Replace across the run boundary naturally merges the
#    runs into one and drops the bookmark that used to sit between
#    them (it gets re-created later, in its new location).
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("This is git test", $false, $false, $false, $false, `
                                 $false, $true, 1, $false, "This is git test", 2)

# ------------------------------------------------------------------
# 2. A new, empty paragraph is inserted right after paragraph 1.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

# ------------------------------------------------------------------
# 3. Another new paragraph is inserted after the empty one; this
#    third paragraph will hold "This is shivani" plus the relocated
#    _GoBack bookmark.
# ------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()

# ------------------------------------------------------------------
# 4. Fill paragraph 3 with the target text. A short placeholder tail
#    is typed along with it so that, below, the bookmark can be
#    inserted at a position that is *not* the very end of the
#    document -- the COM host mis-resolves a collapsed Range sitting
#    exactly at the document's end. The placeholder is stripped back
#    out right after the bookmark is anchored.
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$targetText  = "This is shivani"
$placeholder = "@@TMP@@"
$p3.Range.Text = $targetText + $placeholder

# ------------------------------------------------------------------
# 5. Re-create the _GoBack bookmark, collapsed, right after
#    "This is shivani" (i.e. before the placeholder / paragraph mark).
# ------------------------------------------------------------------
$bmPos = $p2.Range.End + $targetText.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 6. Remove the temporary placeholder text now that the bookmark is
#    safely anchored in place.
# ------------------------------------------------------------------
$tailRange = $d.Range($bmPos, $d.Content.End - 1)
$tailRange.Delete()
